$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 2.02
$ws.Range("R3").Value = 1.88
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.92
$ws.Range("G5").Value = 1.95
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 2.6
$ws.Range("L5").Value = 3.75
$ws.Range("Q5").Value = 1.62
$ws.Range("R5").Value = 2.25
$ws.Range("W5").Value = 10
$ws.Range("Z5").Value = 19
$ws.Range("AA5").Value = 15
$ws.Range("AE5").Value = 12
$ws.Range("AG5").Value = 126
$ws.Range("AJ5").Value = 12
$ws.Range("AN5").Value = 4.33
$ws.Range("AU5").Value = 7
$ws.Range("AW5").Value = 5.5
$ws.Range("G6").Value = 3.1
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 2.38
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 2.25
$ws.Range("U6").Value = 1.53
$ws.Range("V6").Value = 2.38
$ws.Range("AH6").Value = 11
$ws.Range("AR6").Value = 51
$ws.Range("AY6").Value = 17
$ws.Range("BC6").Value = 351
$ws.Range("Q8").Value = 2.35
$ws.Range("R8").Value = 1.57
$ws.Range("G9").Value = 2.25
$ws.Range("K9").Value = 1.91
$ws.Range("L9").Value = 4.33
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 1.48
$ws.Range("AA9").Value = 21
$ws.Range("AE9").Value = 19
$ws.Range("BB9").Value = 351
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("AG11").Value = 1000
$ws.Range("H12").Value = 3.5
$ws.Range("U12").Value = 1.6
$ws.Range("V12").Value = 2.07
$ws.Range("W12").Value = 13
$ws.Range("AI12").Value = 9.75
$ws.Range("AK12").Value = 16.5
$ws.Range("AN12").Value = 5.6
$ws.Range("AU12").Value = 6.7
$ws.Range("AW12").Value = 3.85
$ws.Range("G13").Value = 3.35
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 1.88
$ws.Range("J13").Value = 3.6
$ws.Range("K13").Value = 2.45
$ws.Range("L13").Value = 2.35
$ws.Range("M13").Value = 1.02
$ws.Range("O13").Value = 1.13
$ws.Range("P13").Value = 5.1
$ws.Range("Q13").Value = 1.42
$ws.Range("R13").Value = 2.65
$ws.Range("S13").Value = 1.25
$ws.Range("T13").Value = 3.6
$ws.Range("U13").Value = 1.42
$ws.Range("V13").Value = 2.65
$ws.Range("W13").Value = 17.5
$ws.Range("X13").Value = 24
$ws.Range("Y13").Value = 12
$ws.Range("Z13").Value = 50
$ws.Range("AA13").Value = 24
$ws.Range("AB13").Value = 23
$ws.Range("AD13").Value = 8.75
$ws.Range("AF13").Value = 32
$ws.Range("AG13").Value = 150
$ws.Range("AI13").Value = 12.5
$ws.Range("AJ13").Value = 8.75
$ws.Range("AK13").Value = 18
$ws.Range("AL13").Value = 13
$ws.Range("AM13").Value = 17
$ws.Range("AN13").Value = 5.9
$ws.Range("AO13").Value = 16.5
$ws.Range("AP13").Value = 18.5
$ws.Range("AQ13").Value = 70
$ws.Range("AR13").Value = 75
$ws.Range("AT13").Value = 3.6
$ws.Range("H14").Value = 4.3
$ws.Range("K14").Value = 2.37
$ws.Range("L14").Value = 5.9
$ws.Range("O14").Value = 1.18
$ws.Range("R14").Value = 2.3
$ws.Range("S14").Value = 1.31
$ws.Range("T14").Value = 3.15
$ws.Range("W14").Value = 8.25
$ws.Range("AP14").Value = 15
$ws.Range("AQ14").Value = 19
$ws.Range("AR14").Value = 45
$ws.Range("AT14").Value = 3.15
$ws.Range("AV14").Value = 60
$ws.Range("G15").Value = 1.55
$ws.Range("H15").Value = 4.4
$ws.Range("I15").Value = 4.85
$ws.Range("J15").Value = 1.98
$ws.Range("K15").Value = 2.6
$ws.Range("P15").Value = 5.9
$ws.Range("Q15").Value = 1.32
$ws.Range("R15").Value = 3.1
$ws.Range("U15").Value = 1.4
$ws.Range("V15").Value = 2.75
$ws.Range("Z15").Value = 14.5
$ws.Range("AA15").Value = 10.75
$ws.Range("AC15").Value = 28
$ws.Range("AD15").Value = 10.25
$ws.Range("AE15").Value = 12.5
$ws.Range("AF15").Value = 32
$ws.Range("AG15").Value = 150
$ws.Range("AH15").Value = 26
$ws.Range("AJ15").Value = 16.5
$ws.Range("AN15").Value = 4.1
$ws.Range("AO15").Value = 7.3
$ws.Range("AP15").Value = 11.5
$ws.Range("AQ15").Value = 18.5
$ws.Range("AS15").Value = 90
$ws.Range("AU15").Value = 6.4
$ws.Range("AV15").Value = 35
$ws.Range("BC15").Value = 400
$ws.Range("G16").Value = 3.2
$ws.Range("I16").Value = 2.4
$ws.Range("J16").Value = 3.75
$ws.Range("L16").Value = 3.1
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 8
$ws.Range("O16").Value = 1.4
$ws.Range("P16").Value = 2.75
$ws.Range("Q16").Value = 2.3
$ws.Range("R16").Value = 1.6
$ws.Range("AJ16").Value = 9.5
$ws.Range("AN16").Value = 5
$ws.Range("AX16").Value = 13
$ws.Range("AZ16").Value = 41
$ws.Range("BD16").Value = 126
$ws.Range("BD17").Value = 151
$ws.Range("G18").Value = 1.95
$ws.Range("I18").Value = 3.6
$ws.Range("W18").Value = 8.5
$ws.Range("G19").Value = 1.52
$ws.Range("H19").Value = 3.8
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 2.07
$ws.Range("L19").Value = 5.9
$ws.Range("M19").Value = 1.02
$ws.Range("N19").Value = 8.949999999999999
$ws.Range("O19").Value = 1.27
$ws.Range("P19").Value = 3.1
$ws.Range("Q19").Value = 1.87
$ws.Range("R19").Value = 1.85
$ws.Range("S19").Value = 1.4
$ws.Range("U19").Value = 1.9
$ws.Range("V19").Value = 1.72
$ws.Range("W19").Value = 6.3
$ws.Range("X19").Value = 6.8
$ws.Range("Y19").Value = 8
$ws.Range("AA19").Value = 12.5
$ws.Range("AB19").Value = 28
$ws.Range("AC19").Value = 10
$ws.Range("AD19").Value = 7.5
$ws.Range("AE19").Value = 18
$ws.Range("AF19").Value = 90
$ws.Range("AG19").Value = 800
$ws.Range("AH19").Value = 15
$ws.Range("AJ19").Value = 19
$ws.Range("AL19").Value = 70
$ws.Range("AM19").Value = 65
$ws.Range("AN19").Value = 3.25
$ws.Range("AO19").Value = 7.3
$ws.Range("AQ19").Value = 23
$ws.Range("AR19").Value = 60
$ws.Range("AU19").Value = 8
$ws.Range("AW19").Value = 7.3
$ws.Range("AY19").Value = 40
$ws.Range("AZ19").Value = 250
$ws.Range("AV13").Value = 37
$ws.Range("AW13").Value = 4.3
$ws.Range("AX13").Value = 9
$ws.Range("AY13").Value = 13.5
$ws.Range("AZ13").Value = 28
$ws.Range("BA13").Value = 40
